$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 4: column headers ---
# G-O block gains a "Val Accuracy (mean)" header in O4 (replacing old MAX text)
$ws.Range("O4").Value = "Val Accuracy (mean)"

# Q-W block (Test Set Results) now mirrors the same 6 metric headers as the other block
$ws.Range("R4").Value = "Loss (global)"
$ws.Range("S4").Value = "Accuracy (global)"
$ws.Range("T4").Value = "Accuracy (mean)"
$ws.Range("U4").Value = "Accuracy [0]"
$ws.Range("V4").Value = "Accuracy [1]"
$ws.Range("W4").Value = "Accuracy [2]"

# --- Row 6: Val accuracy min value updated ---
$ws.Range("Y6").Value = 0.5

# --- Row 9 (Swedish): fill in newly-obtained LSTM results ---
$ws.Range("C9").Value = "16h 12m (972m)"
$ws.Range("D9").Value = "1h 01m"
$ws.Range("E9").Value = "1m 13s"

# --- Section headers (merged cells B3:E3, G3:O3, Q3:W3) ---
$ws.Range("G3").Value = "[LSTM] Best Hyperparameters, Best Epoch and Val Accuracy"
$ws.Range("Q3").Value = "[LSTM] Test Set Results"
$ws.Range("B3").Value = "[LSTM] Training Times"

# --- Row 9 (Swedish): remaining hyperparameter / test columns ---
$ws.Range("J9").Value = 1024
$ws.Range("K9").Value = 1024
$ws.Range("L9").Value = 0.001
$ws.Range("M9").Value = 0
$ws.Range("N9").Value = 48
$ws.Range("O9").Value = 0.8459

$ws.Range("R9").Value = 0.9211
$ws.Range("S9").Value = 0.8192
$ws.Range("T9").Value = 0.8196
$ws.Range("U9").Value = 0.8124
$ws.Range("V9").Value = 0.85
$ws.Range("W9").Value = 0.79

# --- Sheet view: scroll/selection position moved ---
$ws.Range("W20").Select()
$excel.ActiveWindow.ScrollColumn = 18
$excel.ActiveWindow.ScrollRow = 1

# --- Workbook window position moved (best effort) ---
$win = $wb.Windows.Item(1)
$win.Left = 28800
$win.Top = 8220
